# Automatische test-sync: 2025-08-05 19:26:50
#
# Appends the new "Testmail #3: Los jij dit even af?" log entry as row 44
# on the "Logs" sheet (extending the conditional formatting ranges to
# match), and refreshes the category counts / ordering on the "Dashboard"
# sheet to reflect that the new mail was categorised as "Overig".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new row of data (row 44)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 44

$logs.Cells.Item($newRow, 1).Value2  = "Los jij dit even af?"
$logs.Cells.Item($newRow, 2).Value2  = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value2  = "Testmail #3: Los jij dit even af?"
$logs.Cells.Item($newRow, 4).Value2  = "Overig"
$logs.Cells.Item($newRow, 5).Value2  = "Beste klant,`nBedankt voor uw e-mail. Het lijkt erop dat u een testmail hebt gestuurd. Als u hulp nodig heeft bij iets specifieks, laat het ons dan weten, zodat we u op de juiste manier kunnen assisteren.`nMet vriendelijke groet,`n[Naam van het bedrijf]"
$logs.Cells.Item($newRow, 6).Value2  = "2025-08-05 19:25:51"
$logs.Cells.Item($newRow, 7).Value2  = "Ja"
$logs.Cells.Item($newRow, 8).Value2  = "Nee"
$logs.Cells.Item($newRow, 9).Value2  = "Ja"
$logs.Cells.Item($newRow, 10).Value2 = "Nee"

# The cell in column E wraps across several lines; let Excel recompute the
# row height back to the sheet's normal auto height instead of leaving a
# one-off custom row height behind.
$logs.Rows.Item($newRow).AutoFit()

# Keep the existing conditional-formatting rules (same priorities, dxfIds,
# formulas) but stretch their range one row further down, just like Excel
# does automatically when a formatted column grows.
function Extend-ConditionalFormatting($col) {
  $oldRange = $logs.Range($col + "2:" + $col + "43")
  $newRange = $logs.Range($col + "2:" + $col + "44")
  $fcs = $oldRange.FormatConditions
  $count = $fcs.Count
  for ($i = 1; $i -le $count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($newRange)
  }
}

Extend-ConditionalFormatting "D"
Extend-ConditionalFormatting "G"
Extend-ConditionalFormatting "H"
Extend-ConditionalFormatting "I"
Extend-ConditionalFormatting "J"

# ---------------------------------------------------------------------
# 2. Dashboard sheet: refresh the "Categorie" / "Aantal" summary table.
#    "Overig" now leads with 3, followed by "Klacht / Probleem" and
#    "Opvolging / Status" (both still 2, but swapped order).
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(6, 1).Value2 = "Overig"
$dash.Cells.Item(6, 2).Value2 = 3

$dash.Cells.Item(7, 1).Value2 = "Klacht / Probleem"
$dash.Cells.Item(7, 2).Value2 = 2

$dash.Cells.Item(8, 1).Value2 = "Opvolging / Status"
$dash.Cells.Item(8, 2).Value2 = 2
